{"js": "// Update the cover-page title block (Title / Subtitle / Author / Date\n// styled paragraphs) to the new lernOS Change Management Leitfaden text.\n//\n// Note: the author's name (\"Simon D\u00fcckert\") also appears later in the\n// document's acknowledgements section, so we must NOT do a blind\n// document-wide text replace. Each of the four styles below is used\n// exactly once (on the cover page), so matching by paragraph style and\n// replacing only the first hit per style is safe and precise.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"lernOS Change Management Leitfaden\",\n  \"Subtitle\": \"TBD\",\n  \"Author\": \"Thomas Jenewein et. al.\",\n  \"Date\": \"Version TBD (TBD)\"\n};\n\nconst applied = {};\n\nfor (const p of paragraphs.items) {\n  const style = p.style;\n  if (Object.prototype.hasOwnProperty.call(replacements, style) && !applied[style]) {\n    p.insertText(replacements[style], \"Replace\");\n    applied[style] = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update title block text on the cover page:\n#   Title    : \"lernOS Template Leitfaden\"                     -> \"lernOS Change Management Leitfaden\"\n#   Subtitle : \"Dein eigener lernOS Leitfaden in einem Sprint\"  -> \"TBD\"\n#   Author   : \"Simon Dueckert\"                                 -> \"Thomas Jenewein et. al.\"\n#   Date     : \"Version 0.2 (01.01.2022)\"                       -> \"Version TBD (TBD)\"\n#\n# The same literal text (the author's name) also shows up later in the\n# document's acknowledgements section, so a document-wide Find/Replace\n# would over-match. Instead, target only the four styled paragraphs on\n# the cover page by style name, each exactly once.\n\n$d = $word.ActiveDocument\n\n$newText = @{\n    \"Title\"    = \"lernOS Change Management Leitfaden\"\n    \"Subtitle\" = \"TBD\"\n    \"Author\"   = \"Thomas Jenewein et. al.\"\n    \"Date\"     = \"Version TBD (TBD)\"\n}\n\n$done = @{}\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($newText.ContainsKey($styleName) -and -not $done.ContainsKey($styleName)) {\n        $rng = $p.Range\n        # Trim the trailing paragraph mark from the range so only the\n        # run text is replaced, leaving the paragraph itself intact.\n        $rng.MoveEnd(1, -1) | Out-Null\n        $rng.Text = $newText[$styleName]\n        $done[$styleName] = $true\n    }\n}\n"}
